# Adds a new "as of" forecast column (AD, dated 2020-05-01) and a new
# target-date row (42, dated 2020-05-15) to both the "cases" and "deaths"
# forecast tables, and fills in the now-observed value for 2020-05-01
# (row 28, column B "Observed").

function Set-TextValue($sheet, $row, $col, $text) {
    # Force the cell to be written as literal text (shared string) instead
    # of letting Excel auto-detect the date-like string and convert it to
    # a date serial number. Re-applying the "Normal" style afterwards
    # drops the temporary text number-format again so the cell is left
    # with no style override, matching the rest of the sheet.
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

function Clear-RangeCells($sheet, $row1, $col1, $row2, $col2) {
    # Materialize empty <c/> cell elements across the given rectangle
    # (mirrors the fully "rectangular" sheet layout used by this workbook,
    # where every row carries a cell for every column even when blank).
    $rng = $sheet.Range($sheet.Cells.Item($row1, $col1), $sheet.Cells.Item($row2, $col2))
    $rng.Style = "Normal"
}

function Update-ForecastSheet($sheet, $observedB28, $adValues, $ad42Value) {
    $lastDataRow = 41        # last pre-existing data row
    $newRow = 42              # new target-date row (2020-05-15)
    $adCol = 30               # column AD
    $firstAdValueRow = 29     # AD29 is the first row in the new column that carries a forecast

    # 1) New column header AD1 = "2020-05-01"
    Set-TextValue $sheet 1 $adCol "2020-05-01"

    # 2) New column AD, rows 2..41: start all as blank cells ...
    Clear-RangeCells $sheet 2 $adCol $lastDataRow $adCol

    # ... then fill in the actual forecast values for rows 29..41
    for ($i = 0; $i -lt $adValues.Length; $i++) {
        $sheet.Cells.Item($firstAdValueRow + $i, $adCol).Value = $adValues[$i]
    }

    # 3) Row 28 ("2020-05-01") is no longer in the future: fill the
    #    "Observed" column (B).
    $sheet.Cells.Item(28, 2).Value = $observedB28

    # 4) New row 42 ("2020-05-15")
    Set-TextValue $sheet $newRow 1 "2020-05-15"
    Clear-RangeCells $sheet $newRow 2 $newRow 29   # B42:AC42 blank
    $sheet.Cells.Item($newRow, $adCol).Value = $ad42Value
}

$wb = $excel.ActiveWorkbook

$casesAD = @(99785, 108368, 117837, 127880, 138725, 148816, 158338, 166448, 174976, 183267, 191803, 198645, 205377)
$deathsAD = @(7020, 7414, 8235, 8802, 9331, 9851, 10301, 10626, 11067, 11454, 11752, 12095, 12379)

$wsCases = $wb.Worksheets.Item("cases")
Update-ForecastSheet $wsCases 91589 $casesAD 211351

$wsDeaths = $wb.Worksheets.Item("deaths")
Update-ForecastSheet $wsDeaths 6329 $deathsAD 12806
